$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the leading date line currently reads
#   "09/0" + "5" + "/2017" + " " + "updated 07/24/2018 " + "– on " + ...
# (five separate runs for the old date/"updated" text, then a run for
# "– on "). It should collapse to a single run "07/24/2018 " immediately
# followed by the unchanged "– on " run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("09/05/2017 updated 07/24/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "07/24/2018", 2) | Out-Null

# The Find/Replace above collapses the old run group down to one run, but
# because the new run and the following "– on " run both carry no explicit
# character formatting, the engine also coalesces them together into a
# single run. Re-introduce the run boundary by dropping a bookmark right on
# it and immediately deleting the bookmark again -- the bookmark forces the
# two sides to stay in distinct runs even after it is removed.
$full = $d.Content.Text
$boundary = $full.IndexOf([char]0x2013 + " on")
if ($boundary -lt 0) { throw "could not locate '- on' boundary after date replace" }
$splitRange = $d.Range($boundary, $boundary)
$d.Bookmarks.Add("TempSplit1", $splitRange) | Out-Null
$d.Bookmarks("TempSplit1").Delete()

# ---------------------------------------------------------------------------
# Change 2: "Maintain the installation u..." -> "Maintain coldbox
# installation u...", keeping the bold formatting throughout, with
# "Maintain ", "coldbox" and " installation u..." ending up as separate
# runs (a collapsed bookmark sits right after "coldbox").
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("Maintain the installation")
if ($idx -lt 0) { throw "could not locate 'Maintain the installation'" }
$theStart = $idx + "Maintain ".Length
$theEnd = $theStart + "the".Length
$theRange = $d.Range($theStart, $theEnd)
$theRange.Text = "coldbox"

# As above, replacing "the" with "coldbox" leaves "Maintain " and "coldbox"
# coalesced into one run since both are bold with no other distinguishing
# formatting. Split them apart the same way.
$full = $d.Content.Text
$mIdx = $full.IndexOf("Maintain ")
if ($mIdx -lt 0) { throw "could not locate 'Maintain ' after coldbox replace" }
$afterMaintain = $mIdx + "Maintain ".Length
$splitRange2 = $d.Range($afterMaintain, $afterMaintain)
$d.Bookmarks.Add("TempSplit2", $splitRange2) | Out-Null
$d.Bookmarks("TempSplit2").Delete()

# ---------------------------------------------------------------------------
# Change 3: the "_GoBack" bookmark moves from right after the inline
# picture (later in the document) to a collapsed position right after the
# new "coldbox" text. Word keeps bookmark names unique, so (re-)adding a
# bookmark named "_GoBack" here automatically removes the old
# bookmarkStart/bookmarkEnd pair that used to sit after the <w:drawing>.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$cIdx = $full.IndexOf("coldbox")
if ($cIdx -lt 0) { throw "could not locate 'coldbox'" }
$afterColdbox = $cIdx + "coldbox".Length
$goBackRange = $d.Range($afterColdbox, $afterColdbox)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

Write-Output "Paragraph 1: $($d.Paragraphs(1).Range.Text)"
Write-Output "Paragraph 2: $($d.Paragraphs(2).Range.Text)"
Write-Output "Bookmark count: $($d.Bookmarks.Count)"
